$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells G1 and H1, matching the style of the existing header row
$ws.Range("G1").Value = "Elapsed Time"
$ws.Range("H1").Value = "CPU"
$ws.Range("A1").Copy()
$ws.Range("G1:H1").PasteSpecial(-4122)

# Updated metric values on row 2
$ws.Range("B2").Value = 0.4124328550386086
$ws.Range("C2").Value = 0.9918673887695629
$ws.Range("D2").Value = 0.4947339447595301

# Updated model description text (dropped n_estimators=50)
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=3))])"

# New Elapsed Time / CPU values
$ws.Range("G2").Value = 0.1218615918667638
$ws.Range("H2").Value = 0.991
